# Ads1299_defRegs.xlsx edit: remove MPU-related register bit settings
# (author's note: "No MPU in code (in another repo), plan to merge with master.")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Row 2 (register @ 45) ---
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1

# --- Row 3 (register @ 46) ---
# G3 gets restyled to match its siblings (D3/E3/H3 = "Neutral" style, index 6)
$ws.Range("D3").Copy()
$ws.Range("G3").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("K3").Value = 1

# --- Row 5 (register @ 48) ---
$ws.Range("J5").Value = 0

# --- Rows 7, 8, 9 (register @ 49) ---
# Re-style D:K on rows 7-9 to match rows 10-13 ("Explanatory Text" style, index 11)
$ws.Range("D10:K10").Copy()
$ws.Range("D7:K7").PasteSpecial($xlPasteFormats)
$ws.Range("D8:K8").PasteSpecial($xlPasteFormats)
$ws.Range("D9:K9").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("K7").Value = 1

$ws.Range("I8").Value = 1
$ws.Range("K8").Value = 1

$ws.Range("D9").Value = 1
$ws.Range("K9").Value = 1

# --- Rows 10-13 (register @ 49) ---
$ws.Range("G10").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("G13").Value = 0

# --- Row 22 (register @ 58) ---
$ws.Range("F22").Value = 0

# --- Restore selection to D8 (matches author's last-saved cursor position) ---
$null = $ws.Range("D8").Select()
